# Sync attendance_reports: swap the order of recorder names in column G
# ("dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com") wherever
# both recorders are listed together for a session row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.Cells(1, 7).End(-4121).Row  # xlDown fallback isn't needed; use UsedRange instead
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
